# Update ubuntu workstation setup
# plus rename files so they show in order

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New note row above the RAM/Disk table -------------------------------
# A9:C9 get a new bold note about Proxmox swap usage.
$ws.Range("A9").Value = "Proxmox uses 8GB for swap !!! + OS ?"
$ws.Range("A9:C9").Font.Bold = $true

# --- prox1 disk allocation (E11) goes from 4GB to 10GB, also bold --------
$ws.Range("E11").Value = 10
$ws.Range("E11").Font.Bold = $true

# --- Raid 1 usage bumped from 20 to 25 for the three k8s servers ---------
$ws.Range("K14").Value = 25
$ws.Range("K15").Value = 25
$ws.Range("K16").Value = 25

# --- Kubernetes worker disk allocation reduced from 8GB to 6GB -----------
$ws.Range("E17").Value = 6
$ws.Range("E18").Value = 6
$ws.Range("E19").Value = 6

# --- Move the active selection to match the author's final cursor spot ---
$ws.Range("P25").Select() | Out-Null
